$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D14").Value = "TransitionToServeState [event] team serving & athlete serving"
$ws.Range("D15").Value = "TransitionToServeState [event] team serving & !athlete serving"
$ws.Range("D16").Value = "TransitionToServeState [event] !team serving"

$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("D17").Select()
